$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws1 = $wb.Worksheets.Item('LP1912')
$ws1.Cells.Item(2,1).Value = 'Última actualización: 08:57:13'
$ws1.Cells.Item(3,1).Value = 'Total filas: 134'
$ws1.Cells.Item(47,1).Value = '05:49:40'
$ws1.Cells.Item(47,2).Value = '07:04'
$ws1.Cells.Item(47,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(47,4).Value = 75
$ws1.Cells.Item(47,5).Value = 'LP1912'
$ws1.Cells.Item(48,1).Value = '05:18:56'
$ws1.Cells.Item(48,2).Value = '07:04'
$ws1.Cells.Item(48,3).Value = '15_ABASTO'
$ws1.Cells.Item(48,4).Value = 106
$ws1.Cells.Item(48,5).Value = 'LP1912'
$ws1.Cells.Item(55,1).Value = '06:43:40'
$ws1.Cells.Item(55,2).Value = '07:16'
$ws1.Cells.Item(55,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(55,4).Value = 33
$ws1.Cells.Item(55,5).Value = 'LP1912'
$ws1.Cells.Item(56,1).Value = '06:15:04'
$ws1.Cells.Item(56,2).Value = '07:16'
$ws1.Cells.Item(56,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(56,4).Value = 61
$ws1.Cells.Item(56,5).Value = 'LP1912'
$ws1.Cells.Item(95,1).Value = '08:57:13'
$ws1.Cells.Item(95,2).Value = '08:59'
$ws1.Cells.Item(95,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(95,4).Value = 2
$ws1.Cells.Item(95,5).Value = 'LP1912'
$ws1.Cells.Item(97,1).Value = '08:21:50'
$ws1.Cells.Item(97,2).Value = '09:01'
$ws1.Cells.Item(97,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(97,4).Value = 40
$ws1.Cells.Item(97,5).Value = 'LP1912'
$ws1.Cells.Item(98,1).Value = '08:57:13'
$ws1.Cells.Item(98,2).Value = '09:02'
$ws1.Cells.Item(98,3).Value = '215A_EL PATO'
$ws1.Cells.Item(98,4).Value = 5
$ws1.Cells.Item(98,5).Value = 'LP1912'
$ws1.Cells.Item(99,1).Value = '07:59:28'
$ws1.Cells.Item(99,2).Value = '09:03'
$ws1.Cells.Item(99,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(99,4).Value = 64
$ws1.Cells.Item(99,5).Value = 'LP1912'
$ws1.Cells.Item(100,1).Value = '08:39:44'
$ws1.Cells.Item(100,2).Value = '09:04'
$ws1.Cells.Item(100,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(100,4).Value = 25
$ws1.Cells.Item(100,5).Value = 'LP1912'
$ws1.Cells.Item(101,1).Value = '08:57:13'
$ws1.Cells.Item(101,2).Value = '09:05'
$ws1.Cells.Item(101,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(101,4).Value = 8
$ws1.Cells.Item(101,5).Value = 'LP1912'
$ws1.Cells.Item(102,1).Value = '08:21:50'
$ws1.Cells.Item(102,2).Value = '09:07'
$ws1.Cells.Item(102,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(102,4).Value = 46
$ws1.Cells.Item(102,5).Value = 'LP1912'
$ws1.Cells.Item(103,1).Value = '07:20:40'
$ws1.Cells.Item(103,2).Value = '09:10'
$ws1.Cells.Item(103,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(103,4).Value = 110
$ws1.Cells.Item(103,5).Value = 'LP1912'
$ws1.Cells.Item(104,1).Value = '08:57:13'
$ws1.Cells.Item(104,2).Value = '09:11'
$ws1.Cells.Item(104,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(104,4).Value = 14
$ws1.Cells.Item(104,5).Value = 'LP1912'
$ws1.Cells.Item(105,1).Value = '08:21:50'
$ws1.Cells.Item(105,2).Value = '09:13'
$ws1.Cells.Item(105,3).Value = '10_OLMOS'
$ws1.Cells.Item(105,4).Value = 52
$ws1.Cells.Item(105,5).Value = 'LP1912'
$ws1.Cells.Item(106,1).Value = '07:20:40'
$ws1.Cells.Item(106,2).Value = '09:16'
$ws1.Cells.Item(106,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(106,4).Value = 116
$ws1.Cells.Item(106,5).Value = 'LP1912'
$ws1.Cells.Item(107,1).Value = '08:57:13'
$ws1.Cells.Item(107,2).Value = '09:17'
$ws1.Cells.Item(107,3).Value = '27_EL RETIRO'
$ws1.Cells.Item(107,4).Value = 20
$ws1.Cells.Item(107,5).Value = 'LP1912'
$ws1.Cells.Item(108,1).Value = '08:21:50'
$ws1.Cells.Item(108,2).Value = '09:21'
$ws1.Cells.Item(108,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(108,4).Value = 60
$ws1.Cells.Item(108,5).Value = 'LP1912'
$ws1.Cells.Item(109,1).Value = '07:59:28'
$ws1.Cells.Item(109,2).Value = '09:22'
$ws1.Cells.Item(109,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(109,4).Value = 83
$ws1.Cells.Item(109,5).Value = 'LP1912'
$ws1.Cells.Item(110,1).Value = '07:47:32'
$ws1.Cells.Item(110,2).Value = '09:22'
$ws1.Cells.Item(110,3).Value = '17_ROMERO'
$ws1.Cells.Item(110,4).Value = 95
$ws1.Cells.Item(110,5).Value = 'LP1912'
$ws1.Cells.Item(111,1).Value = '08:57:13'
$ws1.Cells.Item(111,2).Value = '09:23'
$ws1.Cells.Item(111,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(111,4).Value = 26
$ws1.Cells.Item(111,5).Value = 'LP1912'
$ws1.Cells.Item(112,1).Value = '07:47:32'
$ws1.Cells.Item(112,2).Value = '09:23'
$ws1.Cells.Item(112,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(112,4).Value = 96
$ws1.Cells.Item(112,5).Value = 'LP1912'
$ws1.Cells.Item(113,1).Value = '08:57:13'
$ws1.Cells.Item(113,2).Value = '09:24'
$ws1.Cells.Item(113,3).Value = '11_ETCHEVERRY'
$ws1.Cells.Item(113,4).Value = 27
$ws1.Cells.Item(113,5).Value = 'LP1912'
$ws1.Cells.Item(114,1).Value = '08:21:50'
$ws1.Cells.Item(114,2).Value = '09:29'
$ws1.Cells.Item(114,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(114,4).Value = 68
$ws1.Cells.Item(114,5).Value = 'LP1912'
$ws1.Cells.Item(115,1).Value = '07:47:32'
$ws1.Cells.Item(115,2).Value = '09:32'
$ws1.Cells.Item(115,3).Value = '15_ABASTO'
$ws1.Cells.Item(115,4).Value = 105
$ws1.Cells.Item(115,5).Value = 'LP1912'
$ws1.Cells.Item(116,1).Value = '07:47:32'
$ws1.Cells.Item(116,2).Value = '09:33'
$ws1.Cells.Item(116,3).Value = '10_OLMOS'
$ws1.Cells.Item(116,4).Value = 106
$ws1.Cells.Item(116,5).Value = 'LP1912'
$ws1.Cells.Item(117,1).Value = '08:39:44'
$ws1.Cells.Item(117,2).Value = '09:34'
$ws1.Cells.Item(117,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(117,4).Value = 55
$ws1.Cells.Item(117,5).Value = 'LP1912'
$ws1.Cells.Item(118,1).Value = '08:39:44'
$ws1.Cells.Item(118,2).Value = '09:34'
$ws1.Cells.Item(118,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(118,4).Value = 55
$ws1.Cells.Item(118,5).Value = 'LP1912'
$ws1.Cells.Item(119,1).Value = '08:57:13'
$ws1.Cells.Item(119,2).Value = '09:35'
$ws1.Cells.Item(119,3).Value = '23_HERNANDEZ'
$ws1.Cells.Item(119,4).Value = 38
$ws1.Cells.Item(119,5).Value = 'LP1912'
$ws1.Cells.Item(120,1).Value = '08:57:13'
$ws1.Cells.Item(120,2).Value = '09:35'
$ws1.Cells.Item(120,3).Value = '16_SANTA ANA'
$ws1.Cells.Item(120,4).Value = 38
$ws1.Cells.Item(120,5).Value = 'LP1912'
$ws1.Cells.Item(121,1).Value = '08:21:50'
$ws1.Cells.Item(121,2).Value = '09:41'
$ws1.Cells.Item(121,3).Value = '215C_EL PATO'
$ws1.Cells.Item(121,4).Value = 80
$ws1.Cells.Item(121,5).Value = 'LP1912'
$ws1.Cells.Item(122,1).Value = '07:47:32'
$ws1.Cells.Item(122,2).Value = '09:42'
$ws1.Cells.Item(122,3).Value = '215C_EL PATO'
$ws1.Cells.Item(122,4).Value = 115
$ws1.Cells.Item(122,5).Value = 'LP1912'
$ws1.Cells.Item(123,1).Value = '07:47:32'
$ws1.Cells.Item(123,2).Value = '09:43'
$ws1.Cells.Item(123,3).Value = '14_ABASTO'
$ws1.Cells.Item(123,4).Value = 116
$ws1.Cells.Item(123,5).Value = 'LP1912'
$ws1.Cells.Item(124,1).Value = '08:57:13'
$ws1.Cells.Item(124,2).Value = '09:44'
$ws1.Cells.Item(124,3).Value = '14_ABASTO'
$ws1.Cells.Item(124,4).Value = 47
$ws1.Cells.Item(124,5).Value = 'LP1912'
$ws1.Cells.Item(125,1).Value = '08:49:51'
$ws1.Cells.Item(125,2).Value = '09:52'
$ws1.Cells.Item(125,3).Value = '15_ABASTO'
$ws1.Cells.Item(125,4).Value = 63
$ws1.Cells.Item(125,5).Value = 'LP1912'
$ws1.Cells.Item(126,1).Value = '08:49:51'
$ws1.Cells.Item(126,2).Value = '09:53'
$ws1.Cells.Item(126,3).Value = '10_OLMOS'
$ws1.Cells.Item(126,4).Value = 64
$ws1.Cells.Item(126,5).Value = 'LP1912'
$ws1.Cells.Item(127,1).Value = '08:39:44'
$ws1.Cells.Item(127,2).Value = '10:06'
$ws1.Cells.Item(127,3).Value = '10_OLMOS'
$ws1.Cells.Item(127,4).Value = 87
$ws1.Cells.Item(127,5).Value = 'LP1912'
$ws1.Cells.Item(128,1).Value = '08:21:50'
$ws1.Cells.Item(128,2).Value = '10:10'
$ws1.Cells.Item(128,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(128,4).Value = 109
$ws1.Cells.Item(128,5).Value = 'LP1912'
$ws1.Cells.Item(129,1).Value = '08:57:13'
$ws1.Cells.Item(129,2).Value = '10:11'
$ws1.Cells.Item(129,3).Value = '16_P MOR-SANTA ANA'
$ws1.Cells.Item(129,4).Value = 74
$ws1.Cells.Item(129,5).Value = 'LP1912'
$ws1.Cells.Item(130,1).Value = '08:21:50'
$ws1.Cells.Item(130,2).Value = '10:12'
$ws1.Cells.Item(130,3).Value = '15_ABASTO'
$ws1.Cells.Item(130,4).Value = 111
$ws1.Cells.Item(130,5).Value = 'LP1912'
$ws1.Cells.Item(131,1).Value = '08:49:51'
$ws1.Cells.Item(131,2).Value = '10:20'
$ws1.Cells.Item(131,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(131,4).Value = 91
$ws1.Cells.Item(131,5).Value = 'LP1912'
$ws1.Cells.Item(132,1).Value = '08:39:44'
$ws1.Cells.Item(132,2).Value = '10:21'
$ws1.Cells.Item(132,3).Value = '26_HERNANDEZ'
$ws1.Cells.Item(132,4).Value = 102
$ws1.Cells.Item(132,5).Value = 'LP1912'
$ws1.Cells.Item(133,1).Value = '08:39:44'
$ws1.Cells.Item(133,2).Value = '10:22'
$ws1.Cells.Item(133,3).Value = '17_ROMERO'
$ws1.Cells.Item(133,4).Value = 103
$ws1.Cells.Item(133,5).Value = 'LP1912'
$ws1.Cells.Item(134,1).Value = '08:39:44'
$ws1.Cells.Item(134,2).Value = '10:26'
$ws1.Cells.Item(134,3).Value = '215A_EL PATO'
$ws1.Cells.Item(134,4).Value = 107
$ws1.Cells.Item(134,5).Value = 'LP1912'
$ws1.Cells.Item(135,1).Value = '08:57:13'
$ws1.Cells.Item(135,2).Value = '10:27'
$ws1.Cells.Item(135,3).Value = '215A_EL PATO'
$ws1.Cells.Item(135,4).Value = 90
$ws1.Cells.Item(135,5).Value = 'LP1912'
$ws1.Cells.Item(136,1).Value = '08:49:51'
$ws1.Cells.Item(136,2).Value = '10:41'
$ws1.Cells.Item(136,3).Value = '17_ROMERO'
$ws1.Cells.Item(136,4).Value = 112
$ws1.Cells.Item(136,5).Value = 'LP1912'
$ws1.Cells.Item(137,1).Value = '08:57:13'
$ws1.Cells.Item(137,2).Value = '10:42'
$ws1.Cells.Item(137,3).Value = '17_ROMERO'
$ws1.Cells.Item(137,4).Value = 105
$ws1.Cells.Item(137,5).Value = 'LP1912'
$ws1.Cells.Item(138,1).Value = '08:49:51'
$ws1.Cells.Item(138,2).Value = '10:43'
$ws1.Cells.Item(138,3).Value = '14_ABASTO'
$ws1.Cells.Item(138,4).Value = 114
$ws1.Cells.Item(138,5).Value = 'LP1912'
$ws1.Cells.Item(139,1).Value = '08:57:13'
$ws1.Cells.Item(139,2).Value = '10:44'
$ws1.Cells.Item(139,3).Value = '14_ABASTO'
$ws1.Cells.Item(139,4).Value = 107
$ws1.Cells.Item(139,5).Value = 'LP1912'

# ===== Sheet: LP1912-215 =====
$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws2.Cells.Item(2,1).Value = 'Última actualización: 08:57:13'
$ws2.Cells.Item(3,1).Value = 'Total filas: 19'
$ws2.Cells.Item(20,1).Value = '08:57:13'
$ws2.Cells.Item(20,2).Value = '09:02'
$ws2.Cells.Item(20,3).Value = '215A_EL PATO'
$ws2.Cells.Item(20,4).Value = 5
$ws2.Cells.Item(20,5).Value = 'LP1912'
$ws2.Cells.Item(21,1).Value = '08:21:50'
$ws2.Cells.Item(21,2).Value = '09:41'
$ws2.Cells.Item(21,3).Value = '215C_EL PATO'
$ws2.Cells.Item(21,4).Value = 80
$ws2.Cells.Item(21,5).Value = 'LP1912'
$ws2.Cells.Item(22,1).Value = '07:47:32'
$ws2.Cells.Item(22,2).Value = '09:42'
$ws2.Cells.Item(22,3).Value = '215C_EL PATO'
$ws2.Cells.Item(22,4).Value = 115
$ws2.Cells.Item(22,5).Value = 'LP1912'
$ws2.Cells.Item(23,1).Value = '08:39:44'
$ws2.Cells.Item(23,2).Value = '10:26'
$ws2.Cells.Item(23,3).Value = '215A_EL PATO'
$ws2.Cells.Item(23,4).Value = 107
$ws2.Cells.Item(23,5).Value = 'LP1912'
$ws2.Cells.Item(24,1).Value = '08:57:13'
$ws2.Cells.Item(24,2).Value = '10:27'
$ws2.Cells.Item(24,3).Value = '215A_EL PATO'
$ws2.Cells.Item(24,4).Value = 90
$ws2.Cells.Item(24,5).Value = 'LP1912'

# ===== Sheet: 6203-6173 =====
$ws3 = $wb.Worksheets.Item('6203-6173')
$ws3.Cells.Item(2,1).Value = 'Última actualización: 08:57:13'
$ws3.Cells.Item(3,1).Value = 'Total filas: 28'
$ws3.Cells.Item(32,1).Value = '08:57:13'
$ws3.Cells.Item(32,2).Value = '10:03'
$ws3.Cells.Item(32,3).Value = '215B_LP-P MOR-40 Y 115'
$ws3.Cells.Item(32,4).Value = 66
$ws3.Cells.Item(32,5).Value = 'L6173'
$ws3.Cells.Item(33,1).Value = '08:57:13'
$ws3.Cells.Item(33,2).Value = '10:54'
$ws3.Cells.Item(33,3).Value = '215A_LA PLATA'
$ws3.Cells.Item(33,4).Value = 117
$ws3.Cells.Item(33,5).Value = 'L6173'
